# Update the "Latest HO Xliff Generate Date" / handoff/handback timestamp
# cells to reflect the new report generation run, per commit
# "Generate Report for Handback".

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: Latest HO Xliff Generate Date for first file
$wsOverview.Range("G2").Value = "2016-08-31 07:13:27"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime for first file
$wsZhCn.Range("H2").Value = "2016-08-31 07:13:22"
$wsZhCn.Range("K2").Value = "2016-08-31 07:13:46"

# de-de sheet: Correspond Handback DateTime for first file
$wsDeDe.Range("K2").Value = "2016-08-31 07:13:53"
